$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.861952666666666
$ws.Range("H2").Value = 14.585858
$ws.Range("I2").Value = 0.3995648519435639
$ws.Range("J2").Value = 0.3995648519435638
$ws.Range("M2").Value = 57.48524866666667
$ws.Range("N2").Value = 172.455746
$ws.Range("O2").Value = 0.2249897472933344
$ws.Range("P2").Value = 0.2249897472933344
$ws.Range("Q2").Value = 279.4905580488964
$ws.Range("R2").Value = 2515.415022440068
$ws.Range("S2").Value = 0.08989799506608101
$ws.Range("T2").Value = 0.089897995066081

# Row 3
$ws.Range("G3").Value = 4.861952666666666
$ws.Range("H3").Value = 14.585858
$ws.Range("I3").Value = 0.3995648519435639
$ws.Range("J3").Value = 0.3995648519435638
$ws.Range("M3").Value = 72.97955566666666
$ws.Range("O3").Value = 0.2856324390668287
$ws.Range("P3").Value = 0.2856324390668287
$ws.Range("Q3").Value = 354.8231452856984
$ws.Range("R3").Value = 3193.408307571285
$ws.Range("S3").Value = 0.1141286832260164
$ws.Range("T3").Value = 0.1141286832260164

# Row 4
$ws.Range("G4").Value = 4.861952666666666
$ws.Range("H4").Value = 14.585858
$ws.Range("I4").Value = 0.3995648519435639
$ws.Range("J4").Value = 0.3995648519435638
$ws.Range("M4").Value = 50.94830300000001
$ws.Range("N4").Value = 152.844909
$ws.Range("O4").Value = 0.1994049966359642
$ws.Range("P4").Value = 0.1994049966359642
$ws.Range("Q4").Value = 247.7082376329914
$ws.Range("R4").Value = 2229.374138696922
$ws.Range("S4").Value = 0.07967522795765589
$ws.Range("T4").Value = 0.07967522795765589

# Row 5
$ws.Range("G5").Value = 4.861952666666666
$ws.Range("H5").Value = 14.585858
$ws.Range("I5").Value = 0.3995648519435639
$ws.Range("J5").Value = 0.3995648519435638
$ws.Range("M5").Value = 74.08852933333333
$ws.Range("N5").Value = 222.265588
$ws.Range("O5").Value = 0.2899728170038728
$ws.Range("P5").Value = 0.2899728170038728
$ws.Range("Q5").Value = 360.2149227616115
$ws.Range("R5").Value = 3241.934304854503
$ws.Range("S5").Value = 0.1158629456938106
$ws.Range("T5").Value = 0.1158629456938106

# Row 6
$ws.Range("I6").Value = 0.04932556406896855
$ws.Range("J6").Value = 0.04932556406896854
$ws.Range("M6").Value = 57.48524866666667
$ws.Range("N6").Value = 172.455746
$ws.Range("O6").Value = 0.2249897472933344
$ws.Range("P6").Value = 0.2249897472933344
$ws.Range("Q6").Value = 34.50260792623423
$ws.Range("R6").Value = 310.5234713361081
$ws.Range("S6").Value = 0.01109774619497841
$ws.Range("T6").Value = 0.01109774619497841

# Row 7
$ws.Range("I7").Value = 0.04932556406896855
$ws.Range("J7").Value = 0.04932556406896854
$ws.Range("M7").Value = 72.97955566666666
$ws.Range("O7").Value = 0.2856324390668287
$ws.Range("P7").Value = 0.2856324390668287
$ws.Range("Q7").Value = 43.80228065809622
$ws.Range("S7").Value = 0.01408898117336661
$ws.Range("T7").Value = 0.01408898117336661

# Row 8
$ws.Range("I8").Value = 0.04932556406896855
$ws.Range("J8").Value = 0.04932556406896854
$ws.Range("M8").Value = 50.94830300000001
$ws.Range("N8").Value = 152.844909
$ws.Range("O8").Value = 0.1994049966359642
$ws.Range("P8").Value = 0.1994049966359642
$ws.Range("Q8").Value = 30.57913749506468
$ws.Range("R8").Value = 275.2122374555821
$ws.Range("S8").Value = 0.00983576393723971
$ws.Range("T8").Value = 0.009835763937239709

# Row 9
$ws.Range("I9").Value = 0.04932556406896855
$ws.Range("J9").Value = 0.04932556406896854
$ws.Range("M9").Value = 74.08852933333333
$ws.Range("N9").Value = 222.265588
$ws.Range("O9").Value = 0.2899728170038728
$ws.Range("P9").Value = 0.2899728170038728
$ws.Range("Q9").Value = 44.46788591351378
$ws.Range("R9").Value = 400.210973221624
$ws.Range("S9").Value = 0.01430307276338382
$ws.Range("T9").Value = 0.01430307276338382

# Row 10
$ws.Range("G10").Value = 4.206754333333333
$ws.Range("H10").Value = 12.620263
$ws.Range("I10").Value = 0.3457193616641432
$ws.Range("J10").Value = 0.3457193616641432
$ws.Range("M10").Value = 57.48524866666667
$ws.Range("N10").Value = 172.455746
$ws.Range("O10").Value = 0.2249897472933344
$ws.Range("P10").Value = 0.2249897472933344
$ws.Range("Q10").Value = 241.8263189312442
$ws.Range("R10").Value = 2176.436870381198
$ws.Range("S10").Value = 0.07778331181522848
$ws.Range("T10").Value = 0.07778331181522846

# Row 11
$ws.Range("G11").Value = 4.206754333333333
$ws.Range("H11").Value = 12.620263
$ws.Range("I11").Value = 0.3457193616641432
$ws.Range("J11").Value = 0.3457193616641432
$ws.Range("M11").Value = 72.97955566666666
$ws.Range("O11").Value = 0.2856324390668287
$ws.Range("P11").Value = 0.2856324390668287
$ws.Range("Q11").Value = 307.0070620454911
$ws.Range("R11").Value = 2763.06355840942
$ws.Range("S11").Value = 0.09874866450475629
$ws.Range("T11").Value = 0.09874866450475629

# Row 12
$ws.Range("G12").Value = 4.206754333333333
$ws.Range("H12").Value = 12.620263
$ws.Range("I12").Value = 0.3457193616641432
$ws.Range("J12").Value = 0.3457193616641432
$ws.Range("M12").Value = 50.94830300000001
$ws.Range("N12").Value = 152.844909
$ws.Range("O12").Value = 0.1994049966359642
$ws.Range("P12").Value = 0.1994049966359642
$ws.Range("Q12").Value = 214.3269944212297
$ws.Range("R12").Value = 1928.942949791067
$ws.Range("S12").Value = 0.06893816814962618
$ws.Range("T12").Value = 0.06893816814962618

# Row 13
$ws.Range("G13").Value = 4.206754333333333
$ws.Range("H13").Value = 12.620263
$ws.Range("I13").Value = 0.3457193616641432
$ws.Range("J13").Value = 0.3457193616641432
$ws.Range("M13").Value = 74.08852933333333
$ws.Range("N13").Value = 222.265588
$ws.Range("O13").Value = 0.2899728170038728
$ws.Range("P13").Value = 0.2899728170038728
$ws.Range("Q13").Value = 311.6722418232937
$ws.Range("R13").Value = 2805.050176409643
$ws.Range("S13").Value = 0.1002492171945323
$ws.Range("T13").Value = 0.1002492171945323

# Row 14
$ws.Range("G14").Value = 2.499212666666667
$ws.Range("H14").Value = 7.497638
$ws.Range("I14").Value = 0.2053902223233243
$ws.Range("J14").Value = 0.2053902223233243
$ws.Range("M14").Value = 57.48524866666667
$ws.Range("N14").Value = 172.455746
$ws.Range("O14").Value = 0.2249897472933344
$ws.Range("P14").Value = 0.2249897472933344
$ws.Range("Q14").Value = 143.6678616142165
$ws.Range("R14").Value = 1293.010754527948
$ws.Range("S14").Value = 0.04621069421704651
$ws.Range("T14").Value = 0.0462106942170465

# Row 15
$ws.Range("G15").Value = 2.499212666666667
$ws.Range("H15").Value = 7.497638
$ws.Range("I15").Value = 0.2053902223233243
$ws.Range("J15").Value = 0.2053902223233243
$ws.Range("M15").Value = 72.97955566666666
$ws.Range("O15").Value = 0.2856324390668287
$ws.Range("P15").Value = 0.2856324390668287
$ws.Range("Q15").Value = 182.3914299298384
$ws.Range("R15").Value = 1641.522869368546
$ws.Range("S15").Value = 0.05866611016268933
$ws.Range("T15").Value = 0.05866611016268933

# Row 16
$ws.Range("G16").Value = 2.499212666666667
$ws.Range("H16").Value = 7.497638
$ws.Range("I16").Value = 0.2053902223233243
$ws.Range("J16").Value = 0.2053902223233243
$ws.Range("M16").Value = 50.94830300000001
$ws.Range("N16").Value = 152.844909
$ws.Range("O16").Value = 0.1994049966359642
$ws.Range("P16").Value = 0.1994049966359642
$ws.Range("Q16").Value = 127.3306442027714
$ws.Range("R16").Value = 1145.975797824942
$ws.Range("S16").Value = 0.04095583659144243
$ws.Range("T16").Value = 0.04095583659144243

# Row 17
$ws.Range("G17").Value = 2.499212666666667
$ws.Range("H17").Value = 7.497638
$ws.Range("I17").Value = 0.2053902223233243
$ws.Range("J17").Value = 0.2053902223233243
$ws.Range("M17").Value = 74.08852933333333
$ws.Range("N17").Value = 222.265588
$ws.Range("O17").Value = 0.2899728170038728
$ws.Range("P17").Value = 0.2899728170038728
$ws.Range("Q17").Value = 185.1629909645715
$ws.Range("R17").Value = 1666.466918681144
$ws.Range("S17").Value = 0.05955758135214607
$ws.Range("T17").Value = 0.05955758135214607

